# Update cryptos list - apply scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'67.977.88"
$ws.Cells.Item(2, 5).Value = "  +1.58%  "
$ws.Cells.Item(3, 4).Value = "'3.332.61"
$ws.Cells.Item(3, 5).Value = "  +1.44%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.17%  "
$ws.Cells.Item(5, 4).Value = "'582.51"
$ws.Cells.Item(5, 5).Value = "  +1.59%  "
$ws.Cells.Item(6, 4).Value = "'177.12"
$ws.Cells.Item(6, 5).Value = "  +1.17%  "
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
$ws.Cells.Item(8, 5).Value = "  +1.46%  "
$ws.Cells.Item(9, 4).Value = "'3.329.35"
$ws.Cells.Item(9, 5).Value = "  +1.48%  "
$ws.Cells.Item(10, 5).Value = "  +6.30%  "
$ws.Cells.Item(11, 5).Value = "  +1.78%  "
$ws.Cells.Item(12, 4).Value = "'47.22"
$ws.Cells.Item(12, 5).Value = "  +4.07%  "
$ws.Cells.Item(13, 5).Value = "  +2.54%  "
$ws.Cells.Item(14, 4).Value = "'694.18"
$ws.Cells.Item(14, 5).Value = "  +0.38%  "
$ws.Cells.Item(15, 4).Value = "'3.870.67"
$ws.Cells.Item(15, 5).Value = "  +1.51%  "
$ws.Cells.Item(16, 5).Value = "  +1.80%  "
$ws.Cells.Item(17, 4).Value = "'67.965.70"
$ws.Cells.Item(17, 5).Value = "  +1.50%  "
$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "'3.346.11"
$ws.Cells.Item(18, 5).Value = "  +1.76%  "
$ws.Cells.Item(19, 2).Value = "TRON"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(19, 4).Value = "'0.118"
$ws.Cells.Item(19, 5).Value = "  -0.49%  "
$ws.Cells.Item(20, 4).Value = "'17.48"
$ws.Cells.Item(20, 5).Value = "  +1.33%  "
$ws.Cells.Item(21, 4).Value = "'11.06"
$ws.Cells.Item(21, 5).Value = "  +3.23%  "
$ws.Cells.Item(22, 5).Value = "  +1.41%  "
$ws.Cells.Item(23, 5).Value = "  +4.27%  "
$ws.Cells.Item(24, 4).Value = "'17.07"
$ws.Cells.Item(24, 5).Value = "  +0.96%  "
$ws.Cells.Item(25, 4).Value = "'99.73"
$ws.Cells.Item(25, 5).Value = "  +0.63%  "
$ws.Cells.Item(26, 5).Value = "  +1.35%  "
$ws.Cells.Item(27, 5).Value = "  +0.21%  "
$ws.Cells.Item(28, 5).Value = "  +3.97%  "
$ws.Cells.Item(29, 4).Value = "'33.29"
$ws.Cells.Item(29, 5).Value = "  -0.78%  "
$ws.Cells.Item(30, 4).Value = "'8.57"
$ws.Cells.Item(30, 5).Value = "  +2.56%  "
$ws.Cells.Item(31, 4).Value = "'7.09"
$ws.Cells.Item(31, 5).Value = "  +6.06%  "
$ws.Cells.Item(32, 4).Value = "'566.65"
$ws.Cells.Item(32, 5).Value = "  -0.64%  "
$ws.Cells.Item(33, 4).Value = "'11.01"
$ws.Cells.Item(33, 5).Value = "  +1.76%  "
$ws.Cells.Item(34, 5).Value = "  +2.97%  "
$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).Value = "'57.42"
$ws.Cells.Item(35, 5).Value = "  +4.16%  "
$ws.Cells.Item(36, 2).Value = "Dai"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(36, 4).Value = "'0.999"
$ws.Cells.Item(36, 5).Value = "  -0.14%  "
$ws.Cells.Item(37, 4).Value = "'3.694.95"
$ws.Cells.Item(37, 5).Value = "  -4.60%  "
$ws.Cells.Item(38, 4).Value = "'3.40"
$ws.Cells.Item(38, 5).Value = "  +2.84%  "
$ws.Cells.Item(39, 4).Value = "'34.63"
$ws.Cells.Item(39, 5).Value = "  +9.21%  "
$ws.Cells.Item(40, 4).Value = "'0.133"
$ws.Cells.Item(40, 5).Value = "  +4.17%  "
$ws.Cells.Item(41, 5).Value = "  +3.05%  "
$ws.Cells.Item(42, 4).Value = "'3.18"
$ws.Cells.Item(42, 5).Value = "  +7.45%  "
$ws.Cells.Item(43, 2).Value = "ApeXProtocol"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(43, 4).Value = "'3.35"
$ws.Cells.Item(43, 5).Value = "  -1.36%  "
$ws.Cells.Item(44, 2).Value = "TheGraph"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(44, 4).Value = "'0.337"
$ws.Cells.Item(44, 5).Value = "  +3.58%  "
$ws.Cells.Item(45, 2).Value = "PEPE"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(45, 4).Value = "'0.0₃0673"
$ws.Cells.Item(45, 5).Value = "  +0.61%  "
$ws.Cells.Item(46, 4).Value = "'0.0407"
$ws.Cells.Item(46, 5).Value = "  +0.97%  "
$ws.Cells.Item(47, 5).Value = "  +5.51%  "
$ws.Cells.Item(48, 5).Value = "  +1.32%  "
$ws.Cells.Item(49, 5).Value = "  -0.10%  "
$ws.Cells.Item(50, 5).Value = "  -2.68%  "
$ws.Cells.Item(51, 5).Value = "  +0.79%  "

Write-Host "Updated cryptos list"
